$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("M3").Value = 1.07
$ws.Range("O3").Value = 1.41
$ws.Range("P3").Value = 2.7

# Row 4
$ws.Range("M4").Value = 1.05
$ws.Range("O4").Value = 1.33

# Row 5
$ws.Range("M5").Value = 1.03
$ws.Range("O5").Value = 1.25
$ws.Range("Q5").Value = 1.95
$ws.Range("R5").Value = 1.9

# Row 8
$ws.Range("J8").Value = 2.63
$ws.Range("Q8").Value = 1.85
$ws.Range("R8").Value = 2

# Row 9
$ws.Range("I9").Value = 6
$ws.Range("J9").Value = 2.05
$ws.Range("L9").Value = 5.5
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 15
$ws.Range("O9").Value = 1.2
$ws.Range("P9").Value = 4.33
$ws.Range("AH9").Value = 17

# Row 10
$ws.Range("G10").Value = 1.91
$ws.Range("J10").Value = 2.63

# Row 12
$ws.Range("K12").Value = 1.95

# Row 13
$ws.Range("Q13").Value = 1.9
$ws.Range("R13").Value = 1.95
